$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New warehouse order data for rows 2-6 (row 1 is the header, left untouched).
# Column A holds numeric-looking IDs that must remain text, so we toggle the
# number format to Text before assigning, then restore the Normal style so no
# stray style index is left behind on the cell.
$ids = @("1005726603", "1005931802", "1001030032", "1005789101", "1007004101")
$names = @("МС МАРКЕТ ЕООД", "ЕЛАЦИТЕ - МЕД АД", "КОМЕ ООД", "ПАРТИ ФУУД ДЗЗД", "КРАСИ КАН ЕООД")
$volumes = @(9, 2, 1, 1, 1)
$lats = @(42.78312, 42.69358, 42.95998703256908, 42.69613, 42.54143)
$lons = @(23.50552, 24.01901, 23.35085604339838, 24.07431, 23.49765)
# GPS column is built from the original literal text (not by re-stringifying
# the parsed doubles), since interpolation rounds to fewer significant digits
# than the source data.
$coords = @("42.78312,23.50552", "42.69358,24.01901", "42.95998703256908,23.35085604339838", "42.69613,24.07431", "42.54143,23.49765")

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2

    $idCell = $ws.Cells.Item($row, 1)
    $idCell.NumberFormat = "@"
    $idCell.Value = $ids[$i]
    $idCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $volumes[$i]
    $ws.Cells.Item($row, 4).Value = $coords[$i]
    $ws.Cells.Item($row, 5).Value = $lats[$i]
    $ws.Cells.Item($row, 6).Value = $lons[$i]
}

Write-Host "Warehouse orders refreshed for rows 2-6"
